$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, duplicating the current row 2 content
# (secretar / Secretariat-Administrativ / Bucuresti / ... application) so the
# new applicant entry below starts from the same shape of data, then shifts
# the previously-existing rows down to 3..7.
$ws.Rows(2).Copy()
$ws.Rows(2).Insert()

# New row 2 is a fresh (malformed/terminated) application entry - only the
# experience value differs from the row it was cloned from.
$ws.Range("F2").Value2 = "0 - 1 fdghnta,1 - 5 an"

# The engine's Insert() doesn't keep existing hyperlink anchors in sync with
# the row shift, so rebuild the hyperlinks collection from scratch: clear
# everything, then re-add a mailto link per data row (A2:A7), in the same
# row order the rows themselves were shifted, so relationship ids line up
# the same way Excel would have renumbered them.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:beatrice.dobre@asmi.ro")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:beatrice.dobre@asmi.ro")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:beatrice.dobre@asmi.ro")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:beatrice.dobre@asmi.ro")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:beatrice.dobre@asmi.ro")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:beatrice.dobre@asmi.ro")

# Hyperlinks.Add() stamps a brand-new cell style on its target; restore the
# original shared "Hyperlink" style across the whole email column so every
# row keeps pointing at the same style slot instead of a freshly-minted one.
$ws.Range("A7").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# Match the author's final selection (cell clicked right before saving).
$ws.Range("F2").Select()
$excel.CutCopyMode = $false
